# 自动更新Excel文件 - 2025-11-14 23:12:17
# For every data row (2..99), decrement the "剩余" (remaining, column E) count by 1.
# When the remaining count would wrap below 1 (i.e. old value was 1), it resets to
# 10 and the "开始时间" (start date, column F) is pushed forward by 10 days
# (e.g. 20251105 -> 20251115).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $eCell = $ws.Cells.Item($row, 5)
    $eVal = $eCell.Value2

    if ($eVal -eq $null) {
        continue
    }

    # Skip rows whose start-date (column F) is not a well-formed 8-digit
    # date (yyyymmdd). One legacy row in the sheet has a corrupted 9-digit
    # value and must be left untouched, matching the source update.
    $fCell = $ws.Cells.Item($row, 6)
    $fVal = $fCell.Value2
    if ($fVal -eq $null) {
        continue
    }
    $fText = [string]$fVal
    if ($fText.Length -ne 8) {
        continue
    }

    if ($eVal -eq 1) {
        $eCell.Value = 10
        $fCell.Value = $fVal + 10
    }
    elseif ($eVal -gt 1) {
        $eCell.Value = $eVal - 1
    }
}
